$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1 - title: merge A1:I1, change text, bump row height, re-wrap
# ---------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$r1 = $ws.Range("A1")
$r1.Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Kvareli Municipality"
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------
# Row 2 - subtitle: same text, just shrink the row back to default
# ---------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------
# Row 3 - A3 blank header cell font becomes Sylfaen 11 (years stay)
# ---------------------------------------------------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Sylfaen"
$a3.Font.Size = 11

# ---------------------------------------------------------------
# Row 4 - "family with disabilities Persons " + real figures
# ---------------------------------------------------------------
$a4 = $ws.Range("A4")
$a4.Value = "family with disabilities Persons "
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.ThemeColor = 1
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true
$ws.Rows.Item(4).RowHeight = 24.75

$row4vals = @(614, 546, 511, 524, 536, 573, 596, 597)
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $row4vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.Borders.Item(9).LineStyle = -4142
    $cell.Borders.Item(8).LineStyle = -4142
}

# ---------------------------------------------------------------
# Row 5 - unmerge, "disabilities Persons " + real figures
# ---------------------------------------------------------------
$ws.Range("A5:H5").UnMerge()
$a5 = $ws.Range("A5")
$a5.Value = "disabilities Persons "
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.ThemeColor = 1
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$ws.Rows.Item(5).RowHeight = 21

$row5vals = @(683, 616, 578, 596, 609, 667, 690, 708)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $row5vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.Borders.Item(8).LineStyle = -4142
    if ($cols[$i] -eq "I") {
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = 2
        $cell.Borders.Item(9).ColorIndex = 1
    } else {
        $cell.Borders.Item(9).LineStyle = -4142
    }
}

# ---------------------------------------------------------------
# Row 6 - source note (was row A6 note text), merge A6:H6
# ---------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Font.Bold = $false
$a6.Font.Underline = $false
$a6.Font.ColorIndex = 1
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true

$src = $a6.Characters(1, 7)
$src.Font.Bold = $true
$src.Font.Underline = $true

$ws.Range("A6:H6").Merge()
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------
# Column widths: only column A keeps a custom width now
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.8164063
for ($c = 2; $c -le 16; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 8.73
}

# ---------------------------------------------------------------
# Sheet-level formatting defaults & selection
# ---------------------------------------------------------------
$ws.Range("A1:I1").Select()
